$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.991.58"
Set-TextValue $ws.Range("E2") "  -4.53%  "
Set-TextValue $ws.Range("D3") "1.741.18"
Set-TextValue $ws.Range("E3") "  -5.03%  "
Set-TextValue $ws.Range("E4") "  -0.30%  "
Set-TextValue $ws.Range("D5") "226.70"
Set-TextValue $ws.Range("E5") "  -3.94%  "
Set-TextValue $ws.Range("E6") "  -3.96%  "
Set-TextValue $ws.Range("E7") "  -0.23%  "
Set-TextValue $ws.Range("D8") "0.2732"
Set-TextValue $ws.Range("E8") "  -1.67%  "
Set-TextValue $ws.Range("D9") "23.34"
Set-TextValue $ws.Range("E9") "  -1.33%  "
Set-TextValue $ws.Range("D10") "0.06627"
Set-TextValue $ws.Range("E10") "  -5.30%  "
Set-TextValue $ws.Range("D11") "0.07563"
Set-TextValue $ws.Range("E11") "  -0.87%  "
Set-TextValue $ws.Range("D12") "1.750.98"
Set-TextValue $ws.Range("E12") "  -4.58%  "
Set-TextValue $ws.Range("D13") "4.706"
Set-TextValue $ws.Range("E13") "  -1.39%  "
Set-TextValue $ws.Range("D14") "0.6038"
Set-TextValue $ws.Range("E14") "  -4.55%  "
Set-TextValue $ws.Range("D15") "1.978.51"
Set-TextValue $ws.Range("E15") "  -4.99%  "
Set-TextValue $ws.Range("D16") "74.66"
Set-TextValue $ws.Range("E16") "  -4.35%  "
Set-TextValue $ws.Range("D17") "0.000008720"
Set-TextValue $ws.Range("E17") "  -11.84%  "
Set-TextValue $ws.Range("D18") "27.988.24"
Set-TextValue $ws.Range("E18") "  -3.46%  "
Set-TextValue $ws.Range("E19") "  -5.01%  "
Set-TextValue $ws.Range("E20") "  -0.37%  "
Set-TextValue $ws.Range("D21") "205.74"
Set-TextValue $ws.Range("E21") "  -5.84%  "
Set-TextValue $ws.Range("E22") "  -2.70%  "
Set-TextValue $ws.Range("D23") "6.633"
Set-TextValue $ws.Range("E23") "  -4.32%  "
Set-TextValue $ws.Range("D24") "1.002"
Set-TextValue $ws.Range("E24") "  -0.22%  "
Set-TextValue $ws.Range("D25") "150.26"
Set-TextValue $ws.Range("E25") "  -3.96%  "
Set-TextValue $ws.Range("D26") "8.107"
Set-TextValue $ws.Range("E26") "  +1.20%  "
Set-TextValue $ws.Range("E27") "  -4.69%  "
Set-TextValue $ws.Range("E28") "  -2.55%  "
Set-TextValue $ws.Range("D29") "1.379"
Set-TextValue $ws.Range("E29") "  -3.33%  "
Set-TextValue $ws.Range("E30") "  -5.02%  "
Set-TextValue $ws.Range("D31") "1.394"
Set-TextValue $ws.Range("E31") "  -3.68%  "
Set-TextValue $ws.Range("E32") "  -2.86%  "
Set-TextValue $ws.Range("D33") "3.727"
Set-TextValue $ws.Range("E33") "  -2.08%  "
Set-TextValue $ws.Range("D34") "1.669"
Set-TextValue $ws.Range("E34") "  -4.00%  "
Set-TextValue $ws.Range("E35") "  -5.64%  "
Set-TextValue $ws.Range("D36") "0.6417"
Set-TextValue $ws.Range("E36") "  -1.47%  "
Set-TextValue $ws.Range("E37") "  -4.96%  "
Set-TextValue $ws.Range("D38") "2.722"
Set-TextValue $ws.Range("E38") "  -1.38%  "
Set-TextValue $ws.Range("E39") "  -5.12%  "
Set-TextValue $ws.Range("D40") "1.132.58"
Set-TextValue $ws.Range("E40") "  -1.49%  "
Set-TextValue $ws.Range("D41") "6.185"
Set-TextValue $ws.Range("E41") "  -6.00%  "
Set-TextValue $ws.Range("D42") "0.8777"
Set-TextValue $ws.Range("E42") "  -1.90%  "
Set-TextValue $ws.Range("E43") "  -0.19%  "
Set-TextValue $ws.Range("D44") "99.62"
Set-TextValue $ws.Range("E44") "  -1.48%  "
Set-TextValue $ws.Range("D45") "1.891.83"
Set-TextValue $ws.Range("E45") "  -5.15%  "
Set-TextValue $ws.Range("D46") "59.47"
Set-TextValue $ws.Range("E46") "  -4.76%  "
Set-TextValue $ws.Range("D47") "1.583"
Set-TextValue $ws.Range("E47") "  -2.56%  "
Set-TextValue $ws.Range("D48") "0.00000000108"
Set-TextValue $ws.Range("E48") "  -4.97%  "
Set-TextValue $ws.Range("D49") "8.255"
Set-TextValue $ws.Range("E49") "  -3.25%  "
Set-TextValue $ws.Range("D50") "0.05382"
Set-TextValue $ws.Range("E50") "  -2.10%  "
Set-TextValue $ws.Range("D51") "6.270"
Set-TextValue $ws.Range("E51") "  -2.86%  "
